$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/18/2024  Through  11/24/2024"

# --- Helper: copy number-format/style only (no value) from a donor cell ---
function Copy-Style($donorRef, $targetRef) {
    $ws.Range($donorRef).Copy()
    $ws.Range($targetRef).PasteSpecial(-4122)
}

# --- Helper: force a numeric-looking string (e.g. "0", "***.*") into a cell as
#     TEXT (matching the report convention of showing "0"/"***.*" placeholders),
#     then restore the correct number-format/style from a same-styled donor cell. ---
function Set-TextPlaceholder($targetRef, $donorRef, $text) {
    $ws.Range("ZZ100").Value = "'" + $text
    $ws.Range("ZZ100").Copy()
    $ws.Range($targetRef).PasteSpecial(-4163)
    Copy-Style $donorRef $targetRef
}

# --- Cells that flip from a numeric count/pct to the text placeholders "0"/"***.*" ---
Set-TextPlaceholder "D15" "C15" "0"
Set-TextPlaceholder "E15" "C14" "***.*"
Set-TextPlaceholder "D27" "C27" "0"
Set-TextPlaceholder "E27" "C14" "***.*"
Set-TextPlaceholder "C28" "C27" "0"
Set-TextPlaceholder "G29" "F29" "0"
Set-TextPlaceholder "H29" "C14" "***.*"
Set-TextPlaceholder "G30" "F30" "0"
Set-TextPlaceholder "H30" "C14" "***.*"

# --- Cells that flip from the text placeholders back to real numbers (row 33) ---
Copy-Style "I33" "C33"
$ws.Range("C33").Value = 1
Copy-Style "I33" "D33"
$ws.Range("D33").Value = 1
Copy-Style "K33" "E33"
$ws.Range("E33").Value = 0
Copy-Style "I33" "F33"
$ws.Range("F33").Value = 1
Copy-Style "I33" "G33"
$ws.Range("G33").Value = 1
Copy-Style "K33" "H33"
$ws.Range("H33").Value = 0

# --- Plain numeric value updates ---
$ws.Range("F15").Value = 6
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 21
$ws.Range("K15").Value = -12.5
$ws.Range("L15").Value = 10.526315789473
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = 90.90909090909
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -35
$ws.Range("I16").Value = 215
$ws.Range("J16").Value = 238
$ws.Range("K16").Value = -9.663865546218
$ws.Range("L16").Value = 22.15909090909
$ws.Range("M16").Value = 23.563218390804
$ws.Range("N16").Value = -76.373626373626
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 100
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 235
$ws.Range("J17").Value = 217
$ws.Range("K17").Value = 8.294930875576
$ws.Range("L17").Value = 16.336633663366
$ws.Range("M17").Value = 95.833333333333
$ws.Range("N17").Value = -17.543859649122
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -44.444444444444
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -37.5
$ws.Range("I18").Value = 219
$ws.Range("J18").Value = 201
$ws.Range("K18").Value = 8.955223880597
$ws.Range("L18").Value = 45.033112582781
$ws.Range("M18").Value = -4.366812227074
$ws.Range("N18").Value = -84.345961401
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -7.692307692307
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = 10.63829787234
$ws.Range("I19").Value = 654
$ws.Range("J19").Value = 670
$ws.Range("K19").Value = -2.388059701492
$ws.Range("L19").Value = 5.314009661835
$ws.Range("M19").Value = 55.344418052256
$ws.Range("N19").Value = -21.957040572792
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 236
$ws.Range("J20").Value = 282
$ws.Range("K20").Value = -16.312056737588
$ws.Range("L20").Value = 11.320754716981
$ws.Range("M20").Value = 25.531914893617
$ws.Range("N20").Value = -87.63750654793
$ws.Range("C21").Value = 26
$ws.Range("E21").Value = -25.714285714285
$ws.Range("F21").Value = 126
$ws.Range("G21").Value = 131
$ws.Range("H21").Value = -3.816793893129
$ws.Range("I21").Value = 1581
$ws.Range("J21").Value = 1632
$ws.Range("K21").Value = -3.125
$ws.Range("L21").Value = 14.399421128798
$ws.Range("M21").Value = 37.597911227154
$ws.Range("N21").Value = -70.503731343283
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = -25
$ws.Range("J22").Value = 79
$ws.Range("K22").Value = -31.645569620253
$ws.Range("L22").Value = -28
$ws.Range("M22").Value = 17.391304347826
$ws.Range("C24").Value = 61
$ws.Range("D24").Value = 42
$ws.Range("E24").Value = 45.238095238095
$ws.Range("F24").Value = 163
$ws.Range("G24").Value = 162
$ws.Range("H24").Value = 0.617283950617
$ws.Range("I24").Value = 1855
$ws.Range("J24").Value = 1863
$ws.Range("K24").Value = -0.429414922168
$ws.Range("L24").Value = 34.90909090909
$ws.Range("M24").Value = 117.978848413631
$ws.Range("C25").Value = 37
$ws.Range("D25").Value = 28
$ws.Range("E25").Value = 32.142857142857
$ws.Range("F25").Value = 107
$ws.Range("G25").Value = 96
$ws.Range("H25").Value = 11.458333333333
$ws.Range("I25").Value = 1255
$ws.Range("J25").Value = 1155
$ws.Range("K25").Value = 8.658008658008
$ws.Range("L25").Value = 122.12389380531
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = -16.666666666666
$ws.Range("I26").Value = 507
$ws.Range("J26").Value = 471
$ws.Range("K26").Value = 7.64331210191
$ws.Range("L26").Value = -0.392927308447
$ws.Range("M26").Value = 11.920529801324
$ws.Range("F27").Value = 7
$ws.Range("H27").Value = 250
$ws.Range("I27").Value = 30
$ws.Range("K27").Value = -3.225806451612
$ws.Range("L27").Value = 20
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = -28.571428571428
$ws.Range("I28").Value = 60
$ws.Range("J28").Value = 77
$ws.Range("K28").Value = -22.077922077922
$ws.Range("L28").Value = -27.710843373494
$ws.Range("L29").Value = -16.666666666666
$ws.Range("L30").Value = 0
$ws.Range("I33").Value = 7
$ws.Range("J33").Value = 5
$ws.Range("K33").Value = 40
$ws.Range("L33").Value = -12.5

# --- Clean up scratch cell ---
$ws.Range("ZZ100").Clear()

Write-Host "Edit complete"
